$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hora column: update all data rows from 13 to 14
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "14"
$ws.Range("G2:G51").Style = "Normal"

# Price (D) and Volume(1h) (E) updates
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "310.89"
$ws.Range("E2").Value = "1.26%"
$ws.Range("D2:E2").Style = "Normal"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "37.74"
$ws.Range("E3").Value = "0.56%"
$ws.Range("D3:E3").Style = "Normal"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.114"
$ws.Range("E4").Value = "0.16%"
$ws.Range("D4:E4").Style = "Normal"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07900"
$ws.Range("E5").Value = "0.35%"
$ws.Range("D5:E5").Style = "Normal"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.403"
$ws.Range("E6").Value = "1.38%"
$ws.Range("D6:E6").Style = "Normal"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.904"
$ws.Range("E7").Value = "-3.36%"
$ws.Range("D7:E7").Style = "Normal"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "8.226"
$ws.Range("E8").Value = "0.07%"
$ws.Range("D8:E8").Style = "Normal"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "2.864"
$ws.Range("E9").Value = "-7.56%"
$ws.Range("D9:E9").Style = "Normal"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9258"
$ws.Range("E10").Value = "-0.23%"
$ws.Range("D10:E10").Style = "Normal"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1199"
$ws.Range("E11").Value = "-7.16%"
$ws.Range("D11:E11").Style = "Normal"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1910"
$ws.Range("E12").Value = "1.12%"
$ws.Range("D12:E12").Style = "Normal"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09435"
$ws.Range("E13").Value = "7.95%"
$ws.Range("D13:E13").Style = "Normal"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03367"
$ws.Range("E14").Value = "-1.73%"
$ws.Range("D14:E14").Style = "Normal"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09625"
$ws.Range("E15").Value = "-1.37%"
$ws.Range("D15:E15").Style = "Normal"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001367"
$ws.Range("E16").Value = "-1.83%"
$ws.Range("D16:E16").Style = "Normal"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005881"
$ws.Range("E17").Value = "0.07%"
$ws.Range("D17:E17").Style = "Normal"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.530"
$ws.Range("E18").Value = "-1.44%"
$ws.Range("D18:E18").Style = "Normal"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3431"
$ws.Range("E19").Value = "-0.11%"
$ws.Range("D19:E19").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.82%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.70%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2582"
$ws.Range("E22").Value = "3.44%"
$ws.Range("D22:E22").Style = "Normal"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02099"
$ws.Range("E23").Value = "179.60%"
$ws.Range("D23:E23").Style = "Normal"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04359"
$ws.Range("E24").Value = "0.96%"
$ws.Range("D24:E24").Style = "Normal"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001193"
$ws.Range("E25").Value = "-2.42%"
$ws.Range("D25:E25").Style = "Normal"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004275"
$ws.Range("E26").Value = "-6.94%"
$ws.Range("D26:E26").Style = "Normal"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001297"
$ws.Range("E27").Value = "-63.91%"
$ws.Range("D27:E27").Style = "Normal"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02094"
$ws.Range("E39").Value = "-8.40%"
$ws.Range("D39:E39").Style = "Normal"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05107"
$ws.Range("E40").Value = "1.96%"
$ws.Range("D40:E40").Style = "Normal"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007644"
$ws.Range("E41").Value = "1.79%"
$ws.Range("D41:E41").Style = "Normal"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009102"
$ws.Range("E42").Value = "-7.75%"
$ws.Range("D42:E42").Style = "Normal"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1352"
$ws.Range("E43").Value = "-0.13%"
$ws.Range("D43:E43").Style = "Normal"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002055"
$ws.Range("E44").Value = "-1.98%"
$ws.Range("D44:E44").Style = "Normal"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008573"
$ws.Range("E45").Value = "6.55%"
$ws.Range("D45:E45").Style = "Normal"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006684"
$ws.Range("E46").Value = "2.19%"
$ws.Range("D46:E46").Style = "Normal"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.61%"
$ws.Range("D47:E47").Style = "Normal"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002895"
$ws.Range("E48").Value = "-3.68%"
$ws.Range("D48:E48").Style = "Normal"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001196"
$ws.Range("E49").Value = "-0.52%"
$ws.Range("D49:E49").Style = "Normal"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").Value = "-0.61%"
$ws.Range("D50:E50").Style = "Normal"
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").Value = "-0.61%"
$ws.Range("D51:E51").Style = "Normal"

# Row 48/49: name and link swap (BOLO <-> CoinbaseStockToken)
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
